$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 388, shifting existing rows 388-436 down to 389-437
$ws.Rows.Item(388).Insert()

# Populate the newly inserted row 388 with the new record's data
$ws.Range("A388").Value = 9
$ws.Range("B388").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C388").Value = "Metropolitana"
$ws.Range("D388").Value = 45077
$ws.Range("E388").Value = 13
$ws.Range("F388").Value = 300000001
$ws.Range("G388").Value = "Rabanito"
$ws.Range("H388").Value = "Sin especificar"
$ws.Range("I388").Value = "Primera"
$ws.Range("J388").Value = 7000
$ws.Range("K388").Value = 3000
$ws.Range("L388").Value = 3000
$ws.Range("M388").Value = 3000
$ws.Range("N388").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O388").Value = "Provincia de Chacabuco"
$ws.Range("P388").Value = 30
$ws.Range("Q388").Value = 100
$ws.Range("R388").Value = "Hortaliza"
